$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-5: 2023-09-16 -> 2023-10-05
$ws.Range("C2:C5").Value = 45204
